# Update vm_pu.xlsx results for "case with 380 kV done" commit.
# Rows 2-25 (data rows), columns B-F and I-N get updated voltage magnitude values.
# Column A (bus index), G (slack ref = 1) and H (unused) are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrBF = New-Object 'object[,]' 24,5
$arrBF[0,0] = 1.02
$arrBF[0,1] = 1.048015423470132
$arrBF[0,2] = 1.057568649418278
$arrBF[0,3] = 1.045403463593746
$arrBF[0,4] = 1.065002389984858
$arrBF[1,0] = 1.02
$arrBF[1,1] = 1.049313685056023
$arrBF[1,2] = 1.058810059735478
$arrBF[1,3] = 1.04651783024909
$arrBF[1,4] = 1.066311600093681
$arrBF[2,0] = 1.02
$arrBF[2,1] = 1.050152914912211
$arrBF[2,2] = 1.059612784708318
$arrBF[2,3] = 1.047238488577917
$arrBF[2,4] = 1.067158218374367
$arrBF[3,0] = 1.02
$arrBF[3,1] = 1.05050553150476
$arrBF[3,2] = 1.059950121334191
$arrBF[3,3] = 1.047541357146119
$arrBF[3,4] = 1.067514013297792
$arrBF[4,0] = 1.02
$arrBF[4,1] = 1.050564726013959
$arrBF[4,2] = 1.06000675415845
$arrBF[4,3] = 1.047592204537167
$arrBF[4,4] = 1.067573745683031
$arrBF[5,0] = 1.02
$arrBF[5,1] = 1.050157627356652
$arrBF[5,2] = 1.059617292720002
$arrBF[5,3] = 1.047242535897678
$arrBF[5,4] = 1.067162973004654
$arrBF[6,0] = 1.02
$arrBF[6,1] = 1.04845435153142
$arrBF[6,2] = 1.057988305549838
$arrBF[6,3] = 1.045780155596787
$arrBF[6,4] = 1.065444954880919
$arrBF[7,0] = 1.02
$arrBF[7,1] = 1.045446426628395
$arrBF[7,2] = 1.05511347761806
$arrBF[7,3] = 1.043199999902101
$arrBF[7,4] = 1.062413404458105
$arrBF[8,0] = 1.02
$arrBF[8,1] = 1.043436520855155
$arrBF[8,2] = 1.05319382259821
$arrBF[8,3] = 1.0414775597912
$arrBF[8,4] = 1.060389363028451
$arrBF[9,0] = 1.02
$arrBF[9,1] = 1.042565061419853
$arrBF[9,2] = 1.052361813343167
$arrBF[9,3] = 1.040731136272077
$arrBF[9,4] = 1.059512172916582
$arrBF[10,0] = 1.02
$arrBF[10,1] = 1.042241184209951
$arrBF[10,2] = 1.05205264630125
$arrBF[10,3] = 1.040453788850962
$arrBF[10,4] = 1.059186226309634
$arrBF[11,0] = 1.02
$arrBF[11,1] = 1.042310665116077
$arrBF[11,2] = 1.052118969250421
$arrBF[11,3] = 1.040513285048304
$arrBF[11,4] = 1.059256148418232
$arrBF[12,0] = 1.02
$arrBF[12,1] = 1.042538293293849
$arrBF[12,2] = 1.052336260004254
$arrBF[12,3] = 1.040708212552029
$arrBF[12,4] = 1.059485232521713
$arrBF[13,0] = 1.02
$arrBF[13,1] = 1.042678518753943
$arrBF[13,2] = 1.052470123764449
$arrBF[13,3] = 1.040828301473987
$arrBF[13,4] = 1.059626362892819
$arrBF[14,0] = 1.02
$arrBF[14,1] = 1.043494331983817
$arrBF[14,2] = 1.05324902343968
$arrBF[14,3] = 1.041527084657031
$arrBF[14,4] = 1.06044756275586
$arrBF[15,0] = 1.02
$arrBF[15,1] = 1.044005757109616
$arrBF[15,2] = 1.053737393614402
$arrBF[15,3] = 1.04196525131142
$arrBF[15,4] = 1.060962471664181
$arrBF[16,0] = 1.02
$arrBF[16,1] = 1.044303951411756
$arrBF[16,2] = 1.054022176091502
$arrBF[16,3] = 1.042220769248088
$arrBF[16,4] = 1.061262735600708
$arrBF[17,0] = 1.02
$arrBF[17,1] = 1.044405609217629
$arrBF[17,2] = 1.05411926684693
$arrBF[17,3] = 1.042307884639445
$arrBF[17,4] = 1.061365105457902
$arrBF[18,0] = 1.02
$arrBF[18,1] = 1.043950897582531
$arrBF[18,2] = 1.053685003962379
$arrBF[18,3] = 1.041918246121736
$arrBF[18,4] = 1.060907234466237
$arrBF[19,0] = 1.02
$arrBF[19,1] = 1.042471267429663
$arrBF[19,2] = 1.052272276683894
$arrBF[19,3] = 1.04065081383814
$arrBF[19,4] = 1.059417776283679
$arrBF[20,0] = 1.02
$arrBF[20,1] = 1.041539930936116
$arrBF[20,2] = 1.051383332671253
$arrBF[20,3] = 1.039853392222021
$arrBF[20,4] = 1.058480603516888
$arrBF[21,0] = 1.02
$arrBF[21,1] = 1.042033749597618
$arrBF[21,2] = 1.051854646894611
$arrBF[21,3] = 1.040276172385839
$arrBF[21,4] = 1.058977483379598
$arrBF[22,0] = 1.02
$arrBF[22,1] = 1.043975686585164
$arrBF[22,2] = 1.05370867682445
$arrBF[22,3] = 1.041939485920479
$arrBF[22,4] = 1.06093219400596
$arrBF[23,0] = 1.02
$arrBF[23,1] = 1.046224845225369
$arrBF[23,2] = 1.055857223051351
$arrBF[23,3] = 1.043867432866914
$arrBF[23,4] = 1.063197648645275

$arrIN = New-Object 'object[,]' 24,6
$arrIN[0,0] = 1.041936920771309
$arrIN[0,1] = 1.0530614208751
$arrIN[0,2] = 1.060303261787828
$arrIN[0,3] = 1.048171800965032
$arrIN[0,4] = 1.067716804637489
$arrIN[0,5] = 1.021441961726428
$arrIN[1,0] = 1.042243630574246
$arrIN[1,1] = 1.054006676621121
$arrIN[1,2] = 1.061357448459393
$arrIN[1,3] = 1.049096812363054
$arrIN[1,4] = 1.068840093995745
$arrIN[1,5] = 1.021762703747757
$arrIN[2,0] = 1.042439875891033
$arrIN[2,1] = 1.054617005341316
$arrIN[2,2] = 1.062038466008293
$arrIN[2,3] = 1.049694360622176
$arrIN[2,4] = 1.069565859417402
$arrIN[2,5] = 1.021969629437224
$arrIN[3,0] = 1.042521847636005
$arrIN[3,1] = 1.05487327507446
$arrIN[3,2] = 1.062324502186465
$arrIN[3,3] = 1.049945333618794
$arrIN[3,4] = 1.069870715846305
$arrIN[3,5] = 1.022056474282534
$arrIN[4,0] = 1.04253558000975
$arrIN[4,1] = 1.054916285635223
$arrIN[4,2] = 1.062372513535571
$arrIN[4,3] = 1.049987459295116
$arrIN[4,4] = 1.069921887686152
$arrIN[4,5] = 1.022071047331631
$arrIN[5,0] = 1.042440973282186
$arrIN[5,1] = 1.054620430854901
$arrIN[5,2] = 1.062042289070323
$arrIN[5,3] = 1.049697715061455
$arrIN[5,4] = 1.069569933923961
$arrIN[5,5] = 1.021970790437895
$arrIN[6,0] = 1.042041034217436
$arrIN[6,1] = 1.05338114816463
$arrIN[6,2] = 1.060659760836748
$arrIN[6,3] = 1.048484620415437
$arrIN[6,4] = 1.068096650352124
$arrIN[6,5] = 1.021550485944116
$arrIN[7,0] = 1.041319277211515
$arrIN[7,1] = 1.051187197005056
$arrIN[7,2] = 1.05821494000093
$arrIN[7,3] = 1.046339267865984
$arrIN[7,4] = 1.065492159221481
$arrIN[7,5] = 1.020805106532351
$arrIN[8,0] = 1.040826611915677
$arrIN[8,1] = 1.049717570127654
$arrIN[8,2] = 1.056579093869964
$arrIN[8,3] = 1.044903701172663
$arrIN[8,4] = 1.063750026184927
$arrIN[8,5] = 1.02030494731332
$arrIN[9,0] = 1.040610543277962
$arrIN[9,1] = 1.049079512558034
$arrIN[9,2] = 1.055869302210121
$arrIN[9,3] = 1.044280788876923
$arrIN[9,4] = 1.062994246877359
$arrIN[9,5] = 1.020087593790648
$arrIN[10,0] = 1.040529872713435
$arrIN[10,1] = 1.048842251185383
$arrIN[10,2] = 1.055605431393658
$arrIN[10,3] = 1.044049213029625
$arrIN[10,4] = 1.062713299309183
$arrIN[10,5] = 1.020006740675637
$arrIN[11,0] = 1.040547195534262
$arrIN[11,1] = 1.048893156296928
$arrIN[11,2] = 1.05566204270575
$arrIN[11,3] = 1.044098895862161
$arrIN[11,4] = 1.062773573414438
$arrIN[11,5] = 1.020024089313857
$arrIN[12,0] = 1.040603883457062
$arrIN[12,1] = 1.049059905747654
$arrIN[12,2] = 1.055847495126178
$arrIN[12,3] = 1.044261650810504
$arrIN[12,4] = 1.062971028120706
$arrIN[12,5] = 1.0200809128666
$arrIN[13,0] = 1.040638755987686
$arrIN[13,1] = 1.049162611257623
$arrIN[13,2] = 1.055961728888018
$arrIN[13,3] = 1.044361903120169
$arrIN[13,4] = 1.063092657558078
$arrIN[13,5] = 1.020115908022937
$arrIN[14,0] = 1.040840893821058
$arrIN[14,1] = 1.049759879841278
$arrIN[14,2] = 1.056626169372291
$arrIN[14,3] = 1.044945014130006
$arrIN[14,4] = 1.063800154462949
$arrIN[14,5] = 1.020319355808283
$arrIN[15,0] = 1.04096695477242
$arrIN[15,1] = 1.050134073503079
$arrIN[15,2] = 1.05704256225098
$arrIN[15,3] = 1.045310433583154
$arrIN[15,4] = 1.064243564807592
$arrIN[15,5] = 1.020446763382046
$arrIN[16,0] = 1.041040219577958
$arrIN[16,1] = 1.05035217057334
$arrIN[16,2] = 1.057285296579102
$arrIN[16,3] = 1.045523451092435
$arrIN[16,4] = 1.064502061384032
$arrIN[16,5] = 1.020521002747647
$arrIN[17,0] = 1.041065156169903
$arrIN[17,1] = 1.050426508353206
$arrIN[17,2] = 1.057368038964752
$arrIN[17,3] = 1.045596063359714
$arrIN[17,4] = 1.064590178874484
$arrIN[17,5] = 1.020546303706603
$arrIN[18,0] = 1.040953456984218
$arrIN[18,1] = 1.050093943026393
$arrIN[18,2] = 1.056997901803204
$arrIN[18,3] = 1.045271240533974
$arrIN[18,4] = 1.06419600528938
$arrIN[18,5] = 1.020433101552833
$arrIN[19,0] = 1.040587201688674
$arrIN[19,1] = 1.049010809347832
$arrIN[19,2] = 1.055792890182
$arrIN[19,3] = 1.0442137290366
$arrIN[19,4] = 1.06291288866304
$arrIN[19,5] = 1.020064183021091
$arrIN[20,0] = 1.040354532051892
$arrIN[20,1] = 1.048328303765421
$arrIN[20,2] = 1.055033962048308
$arrIN[20,3] = 1.043547679804787
$arrIN[20,4] = 1.062104882270888
$arrIN[20,5] = 1.019831544144073
$arrIN[21,0] = 1.040478101614547
$arrIN[21,1] = 1.048690255761387
$arrIN[21,2] = 1.055436407554654
$arrIN[21,3] = 1.043900875032743
$arrIN[21,4] = 1.062533342363036
$arrIN[21,5] = 1.019954935681392
$arrIN[22,0] = 1.040959556870666
$arrIN[22,1] = 1.050112076767289
$arrIN[22,2] = 1.057018082371449
$arrIN[22,3] = 1.045288950573454
$arrIN[22,4] = 1.0642174958122
$arrIN[22,5] = 1.020439274978078
$arrIN[23,0] = 1.041507890492645
$arrIN[23,1] = 1.051755607385358
$arrIN[23,2] = 1.0588480241907
$arrIN[23,3] = 1.046894822487582
$arrIN[23,4] = 1.066166493106689
$arrIN[23,5] = 1.020998372547424

$ws.Range("B2:F25").Value = $arrBF
$ws.Range("I2:N25").Value = $arrIN

Write-Host "Updated vm_pu values for case with 380 kV"
